# Generate Report for Handoff
# The file "a33124e5-5ffc-40cc-b6c3-777194ec2b40.md" (row 3 in every sheet) is
# now ready for handoff, so its status / datetime stamps move forward and a
# version-mismatch error message is recorded.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cf6479b8bf5a799dc518ce1972b15add54059681/e2e/a33124e5-5ffc-40cc-b6c3-777194ec2b40.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e31cbf3d1c84a1382fe51aa861c5f82a695ba629/e2e/a33124e5-5ffc-40cc-b6c3-777194ec2b40.md."

# --- Overview sheet ---
$overview = $wb.Sheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-07 07:08:27"

# --- zh-cn sheet ---
$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-09-07 07:08:21"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.15

# --- de-de sheet ---
$dede = $wb.Sheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-09-07 07:08:27"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.15
